# Update header row to reflect renamed/re-derived model-output columns
# (HybridLR replaces the old "logits" naming; Error_Categories replaces
# the old per-sample "Error_UniCategories" column), as part of the
# "update paths to cloud" pass over the model interpretation outputs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 'HybridLR$_{h}$_y_trues'
$ws.Range("F1").Value = 'HybridLR$_{h}$_pos_matches'
$ws.Range("G1").Value = 'HybridLR$_{h}$_neg_matches'
$ws.Range("E1").Value = 'Error_Categories'

$ws.Range("E1").Select()
